$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "326.47"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value2 = "-2.47%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "44.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value2 = "0.36%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "5.557"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value2 = "-3.47%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "0.08013"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value2 = "-4.45%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "4.296"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value2 = "-5.08%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "1.885"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value2 = "-3.49%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "2.601"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value2 = "-8.70%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.9433"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value2 = "-0.46%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.1158"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value2 = "-6.71%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value2 = "-7.60%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.09698"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value2 = "-3.34%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "0.04369"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value2 = "-1.18%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "0.1064"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value2 = "-0.43%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.001272"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value2 = "-2.02%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "0.04219"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value2 = "-4.22%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "0.005975"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value2 = "-1.53%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "3.602"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value2 = "3.03%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "8.612"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value2 = "-0.95%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "0.1379"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value2 = "1.18%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value2 = "0.52%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "0.001250"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value2 = "-0.82%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "0.004498"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value2 = "3.18%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "0.0001262"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value2 = "-0.16%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "0.0003997"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value2 = "0.00%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.02605"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value2 = "-7.88%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.05389"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value2 = "-8.50%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.007601"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value2 = "-4.52%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.1389"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value2 = "-2.68%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.007289"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value2 = "-19.28%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.002019"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value2 = "-5.99%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.008836"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value2 = "-14.81%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "0.00006924"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value2 = "-4.24%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.00000000751"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value2 = "0.00%"
$ws.Range("B47").Value2 = "CoinbaseStockToken"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.002274"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value2 = "0.00%"
$ws.Range("B48").Value2 = "BOLO"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.003635"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value2 = "13.43%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.00002104"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value2 = "0.00%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.0002004"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value2 = "0.00%"
